# "Removed some not needed components"
#
# The BOM lists reference designators for each passive-component row in
# column C. Two resistors (R56 from the 100k row, R58 from the 1k row)
# are no longer populated on the board, so their designators are dropped
# from the corresponding note strings. Downstream formulas (counts,
# extended cost, concatenated Digikey strings, the cost totals) all
# recalculate automatically from that edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 35: 100k resistors -> remove "R56" ------------------------------
# Plain (non rich-text) string: "R11,R14,R17,R20,R35,R37,R38,R48,R49,R55,R56"
$c35 = $ws.Range("C35")
$full35 = $c35.Value2
$i = $full35.IndexOf(",R56")
if ($i -ge 0) {
    $c35.Characters($i + 1, 4).Text = ""
}

# --- Row 30: 1k resistors -> remove "R58" --------------------------------
# Rich-text string with colored runs highlighting R39/R59 (green, hall
# sensor only) and R64 (red). "R58," sits inside a plain black run
# between R39 and R59; deleting the substring collapses all runs to
# plain text, so the colors are re-applied afterwards to match the
# original formatting exactly.
$c30 = $ws.Range("C30")
$full30 = $c30.Value2
$i = $full30.IndexOf("R58,")
if ($i -ge 0) {
    $c30.Characters($i + 1, 4).Text = ""
}

$newfull30 = $c30.Value2

$iR39 = $newfull30.IndexOf("R39")
$c30.Characters($iR39 + 1, 3).Font.Color = 5287936   # RGB(00,B0,50) green

$iR59 = $newfull30.IndexOf("R59")
$midLen = $iR59 - ($iR39 + 3)
$c30.Characters($iR39 + 4, $midLen).Font.Color = 0   # black
$c30.Characters($iR59 + 1, 3).Font.Color = 5287936   # green

$iR64 = $newfull30.IndexOf("R64")
$mid2Len = $iR64 - ($iR59 + 3)
$c30.Characters($iR59 + 4, $mid2Len).Font.Color = 0  # black
$c30.Characters($iR64 + 1, 3).Font.Color = 255       # red

$tailStart = $iR64 + 4
$tailLen = $newfull30.Length - $tailStart + 1
if ($tailLen -gt 0) {
    $c30.Characters($tailStart, $tailLen).Font.Color = 0   # black
}

# --- Keep the active selection in sync with the author's final position --
$ws.Range("C32").Select()
